# Update the fitted parameter values on Sheet1 (column B) to the newly
# curve-fitted numbers referenced in the commit "tried to fit the curves to tests".
# Note: literal scientific notation (e.g. 5E-5) is not supported by the
# script parser here, so plain decimal literals are used instead.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("B2").Value  = 0.57255529999999999    # var_I
$ws.Range("B3").Value  = 0.95                   # var_PR
$ws.Range("B4").Value  = 0.000050000000000000002 # var_c_a
$ws.Range("B5").Value  = 0.000069999999999999994 # var_c_r
$ws.Range("B6").Value  = 0.037749999999999999   # var_h_C
$ws.Range("B8").Value  = 0.00068110000000000002 # var_h_cham
$ws.Range("B11").Value = 0.0018                 # var_r_A
$ws.Range("B12").Value = 0.0074999999999999997  # var_r_Cout
$ws.Range("B13").Value = 0.001                  # var_r_cham
$ws.Range("B15").Value = 0.001                  # var_t_cfr
$ws.Range("B16").Value = 0.0043499999999999997  # var_t_mba
$ws.Range("B17").Value = 0.00165                # var_t_mta

# Reflect the saved window position recorded by Excel for this session.
$win = $excel.ActiveWindow
$win.Left = 330
$win.Top = 345
